# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Updates the "Periodo Mora" (column E), "Valor Mora" (column F) and
# "Salario Basico" (column G) values for rows 16-44 of Hoja1, rolling the
# period window forward (from 2003..1711 down to 1711..2003 ascending) and
# refreshing the mora/salary amounts that go with the new periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New ascending list of periods (Periodo Mora) for rows 16..44.
$periods = @(
    "1711", "1712",
    "1801", "1802", "1803", "1804", "1805", "1806",
    "1807", "1808", "1809", "1810", "1811", "1812",
    "1901", "1902", "1903", "1904", "1905", "1906",
    "1907", "1908", "1909", "1910", "1911", "1912",
    "2001", "2002", "2003"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $period = $periods[$i]

    $ws.Cells.Item($row, 5).Value = $period

    # Valor Mora: 29509 through period 1808, 31249 from period 1809 onward.
    if ([int]$period -le 1808) {
        $ws.Cells.Item($row, 6).Value = 29509
    } else {
        $ws.Cells.Item($row, 6).Value = 31249
    }

    # Salario Basico: updated uniformly to 781242.
    $ws.Cells.Item($row, 7).Value = 781242
}
